$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.034.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.10%  '
$ws.Range("D3").Value = '''2.258.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.51%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''298.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.64%  '
$ws.Range("D6").Value = '''93.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.64%  '
$ws.Range("D7").Value = '''0.497'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.53%  '
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("E10").Value = '  -5.81%  '
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").Value = '''47.65'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.61%  '
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("D15").Value = '''2.609.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.65%  '
$ws.Range("D16").Value = '''15.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.02%  '
$ws.Range("D17").Value = '''2.247.72'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.27%  '
$ws.Range("D18").Value = '''0.774'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.06%  '
$ws.Range("D19").Value = '''42.030.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("D22").Value = '''11.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.11%  '
$ws.Range("D23").Value = '''66.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = '''233.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '''1.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.71%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -4.50%  '
$ws.Range("D28").Value = '''23.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.98%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '''167.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.76%  '
$ws.Range("D31").Value = '''33.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("D32").Value = '''9.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '''4.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.03%  '
$ws.Range("D35").Value = '''4.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.09%  '
$ws.Range("E36").Value = '  -5.80%  '
$ws.Range("D37").Value = '''0.0692'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.01%  '
$ws.Range("D38").Value = '''2.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.75%  '
$ws.Range("E39").Value = '  -8.52%  '
$ws.Range("E40").Value = '  -3.73%  '
$ws.Range("E41").Value = '  -3.29%  '
$ws.Range("D42").Value = '''1.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.71%  '
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = '''1.949.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.48%  '
$ws.Range("E45").Value = '  -2.38%  '
$ws.Range("D46").Value = '''17.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.71%  '
$ws.Range("E47").Value = '  -7.45%  '
$ws.Range("E48").Value = '  -4.65%  '
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("D50").Value = '''2.483.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.96%  '
$ws.Range("D51").Value = '''52.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.97%  '
